$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310.14816
$ws.Range("I33").Value = 295.6087
$ws.Range("J33").Value = 393.75
$ws.Range("K33").Value = 295.6087
$ws.Range("L33").Value = 393.75
$ws.Range("M33").Value = -66.6087
$ws.Range("N33").Value = -851.75

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1310.5
$ws.Range("I80").Value = 2830
$ws.Range("J80").Value = 804
$ws.Range("K80").Value = 8490
$ws.Range("L80").Value = 2412
$ws.Range("M80").Value = -7492
$ws.Range("N80").Value = -4408

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1310.5
$ws.Range("I83").Value = 2830
$ws.Range("J83").Value = 804
$ws.Range("K83").Value = 25470
$ws.Range("L83").Value = 7236
$ws.Range("M83").Value = -20478
$ws.Range("N83").Value = -17220

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1231
$ws.Range("J111").Value = 1857
$ws.Range("L111").Value = 5571
$ws.Range("N111").Value = -11705

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2482.6553
$ws.Range("I113").Value = 2351.25
$ws.Range("J113").Value = 2503.68
$ws.Range("K113").Value = 2351.25
$ws.Range("L113").Value = 2503.68
$ws.Range("M113").Value = 902.75
$ws.Range("N113").Value = -9011.68

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3040.2
$ws.Range("I116").Value = 2307
$ws.Range("J116").Value = 4751
$ws.Range("K116").Value = 2307
$ws.Range("L116").Value = 4751
$ws.Range("M116").Value = 1135
$ws.Range("N116").Value = -11635

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1792.2307
$ws.Range("I2").Value = 933.9091
$ws.Range("J2").Value = 6513
$ws.Range("K2").Value = 933.9091
$ws.Range("L2").Value = 6513
$ws.Range("M2").Value = -820.9091
$ws.Range("N2").Value = -6739

# ARM row 53
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 6000

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1792.2307
$ws.Range("I116").Value = 933.9091
$ws.Range("J116").Value = 6513
$ws.Range("K116").Value = 933.9091
$ws.Range("L116").Value = 6513
$ws.Range("M116").Value = 1360.0909
$ws.Range("N116").Value = -11101

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 27750
$ws.Range("J133").Value = 27750
$ws.Range("L133").Value = 27750
$ws.Range("N133").Value = -32810

# ARM row 140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 31429
$ws.Range("J140").Value = 31429
$ws.Range("L140").Value = 31429
$ws.Range("N140").Value = -41789

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1792.2307
$ws.Range("I3").Value = 933.9091
$ws.Range("J3").Value = 6513
$ws.Range("K3").Value = 933.9091
$ws.Range("L3").Value = 6513
$ws.Range("M3").Value = -819.9091
$ws.Range("N3").Value = -6741

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 768.61536
$ws.Range("I80").Value = 461.33334
$ws.Range("K80").Value = 461.33334
$ws.Range("M80").Value = 536.66666

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 768.61536
$ws.Range("I83").Value = 461.33334
$ws.Range("K83").Value = 2306.6667
$ws.Range("M83").Value = 2685.3333

# BSM row 88
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 34000
$ws.Range("J88").Value = 34000
$ws.Range("L88").Value = 34000
$ws.Range("N88").Value = -34812

# BSM row 91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 34000
$ws.Range("J91").Value = 34000
$ws.Range("L91").Value = 34000
$ws.Range("N91").Value = -36808

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55556624
$ws.Range("I16").Value = 90910010
$ws.Range("J16").Value = 1313
$ws.Range("K16").Value = 90910010
$ws.Range("L16").Value = 1313
$ws.Range("M16").Value = -90909723
$ws.Range("N16").Value = -1887

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2151.3333
$ws.Range("I99").Value = 1906.8572
$ws.Range("K99").Value = 1906.8572
$ws.Range("M99").Value = -408.8571999999999

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 55556624
$ws.Range("I113").Value = 90910010
$ws.Range("J113").Value = 1313
$ws.Range("K113").Value = 90910010
$ws.Range("L113").Value = 1313
$ws.Range("M113").Value = -90907840
$ws.Range("N113").Value = -5653

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2151.3333
$ws.Range("I126").Value = 1906.8572
$ws.Range("K126").Value = 5720.571599999999
$ws.Range("M126").Value = -3250.571599999999

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 685.3333
$ws.Range("I6").Value = 142.66667
$ws.Range("K6").Value = 428.00001
$ws.Range("M6").Value = -315.00001

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1387.5714
$ws.Range("I68").Value = 730.2727
$ws.Range("J68").Value = 1812.8823
$ws.Range("K68").Value = 2190.8181
$ws.Range("L68").Value = 5438.6469
$ws.Range("M68").Value = -1379.8181
$ws.Range("N68").Value = -7060.6469

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1387.5714
$ws.Range("I71").Value = 730.2727
$ws.Range("J71").Value = 1812.8823
$ws.Range("K71").Value = 6572.454299999999
$ws.Range("L71").Value = 16315.9407
$ws.Range("M71").Value = -2516.454299999999
$ws.Range("N71").Value = -24427.9407

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 23810954
$ws.Range("I131").Value = 142857810
$ws.Range("J131").Value = 1585.0857
$ws.Range("K131").Value = 428573430
$ws.Range("L131").Value = 4755.257100000001
$ws.Range("M131").Value = -428568390
$ws.Range("N131").Value = -14835.2571

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1264.7646
$ws.Range("I22").Value = 1192.4615
$ws.Range("J22").Value = 1499.75
$ws.Range("K22").Value = 1192.4615
$ws.Range("L22").Value = 1499.75
$ws.Range("M22").Value = -897.4614999999999
$ws.Range("N22").Value = -2089.75

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1264.7646
$ws.Range("I27").Value = 1192.4615
$ws.Range("J27").Value = 1499.75
$ws.Range("K27").Value = 1192.4615
$ws.Range("L27").Value = 1499.75
$ws.Range("M27").Value = -1085.4615
$ws.Range("N27").Value = -1713.75

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2020.6
$ws.Range("I100").Value = 1925.75
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 1925.75
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -1384.75
$ws.Range("N100").Value = -3482

# LTW row 110
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 15626851
$ws.Range("I122").Value = 19232562
$ws.Range("J122").Value = 2101.6667
$ws.Range("K122").Value = 57697686
$ws.Range("L122").Value = 6305.000100000001
$ws.Range("M122").Value = -57695236
$ws.Range("N122").Value = -11205.0001
